$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.809332370758057
$ws.Range("B1").Value = 2.569871425628662
$ws.Range("C1").Value = 2.098279476165771
$ws.Range("D1").Value = 1.80638325214386
$ws.Range("E1").Value = 1.759797811508179
